# Add a new row of contact data (Isha / J / abcd@xyz.com / 9876543210 /
# 8888888888 / Isha) under the existing "Ranjini" row on Sheet1, matching
# how the rest of the sheet is laid out (First Name, Last Name, Email id,
# Number, Phone Number, Name -> repeats First Name in the last column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Isha"
$ws.Range("B3").Value = "J"
$ws.Range("C3").Value = "abcd@xyz.com"
# The phone-number-like values are stored as text (quote-prefixed), same
# as the existing Number / Phone Number cells in row 2, so the leading
# apostrophe keeps them as literal text rather than numbers.
$ws.Range("D3").Value = "'9876543210"
$ws.Range("E3").Value = "'8888888888"
$ws.Range("F3").Value = "Isha"

# Extend the selection to cover the newly populated range, same as the
# original sheet had the (then) full data range selected.
$ws.Range("A1:F3").Select() | Out-Null
